{"js": "// Remove the trailing \"spacer\" paragraph plus the \"Ver no Jupiter...\" and\n// \"\u00a9 2020 ...\" footer paragraphs that follow the Bibliografia text, while\n// keeping the blank paragraph that comes right after them (and the final\n// page-break paragraph / section properties untouched).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst JUPITER_TEXT = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst COPYRIGHT_TEXT =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n// Locate the \"Ver no Jupiter ...\" paragraph by its exact text.\nlet jupiterIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === JUPITER_TEXT) {\n    jupiterIdx = i;\n    break;\n  }\n}\n\nif (jupiterIdx === -1) {\n  throw new Error(\"Could not find the 'Ver no Jupiter ...' paragraph.\");\n}\n\n// Sanity-check the paragraph right after it is the copyright line, and the\n// paragraph right before it is the empty spacer paragraph to remove too.\nconst copyrightIdx = jupiterIdx + 1;\nif (!items[copyrightIdx] || items[copyrightIdx].text !== COPYRIGHT_TEXT) {\n  throw new Error(\"Could not find the '\u00a9 2020 ...' paragraph after it.\");\n}\n\nconst spacerIdx = jupiterIdx - 1;\nif (!items[spacerIdx] || items[spacerIdx].text !== \"\") {\n  throw new Error(\"Could not find the empty spacer paragraph before it.\");\n}\n\n// Delete the three paragraphs (spacer, Jupiter line, copyright line).\n// The shim retargets pinned paragraph indices across a batch of deletes\n// queued in the same sync, so queuing all three before `sync()` is safe.\nitems[spacerIdx].delete();\nitems[jupiterIdx].delete();\nitems[copyrightIdx].delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"spacer\" paragraph plus the \"Ver no Jupiter...\" and\n# \"\u00a9 2020 ...\" footer paragraphs that follow the Bibliografia text, while\n# keeping the blank paragraph that comes right after them (and the final\n# page-break paragraph / section properties untouched).\n\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n# Locate the \"Ver no Jupiter ...\" paragraph by its exact text (Paragraph\n# Range.Text carries a trailing paragraph mark, so trim it before comparing).\n$jupiterIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\")\n    if ($t -eq $jupiterText) {\n        $jupiterIdx = $i\n        break\n    }\n}\nif ($jupiterIdx -eq -1) {\n    throw \"Could not find the 'Ver no Jupiter ...' paragraph.\"\n}\n\n# Sanity-check the paragraph right after it is the copyright line, and the\n# paragraph right before it is the empty spacer paragraph to remove too.\n$copyrightIdx = $jupiterIdx + 1\n$copyrightActual = $d.Paragraphs.Item($copyrightIdx).Range.Text.TrimEnd(\"`r\")\nif ($copyrightActual -ne $copyrightText) {\n    throw \"Could not find the '\u00a9 2020 ...' paragraph after it.\"\n}\n\n$spacerIdx = $jupiterIdx - 1\n$spacerActual = $d.Paragraphs.Item($spacerIdx).Range.Text.TrimEnd(\"`r\")\nif ($spacerActual -ne \"\") {\n    throw \"Could not find the empty spacer paragraph before it.\"\n}\n\n# Delete the three paragraphs (copyright, Jupiter line, spacer) starting from\n# the highest index so the lower, not-yet-deleted indices stay valid.\n$d.Paragraphs.Item($copyrightIdx).Range.Delete()\n$d.Paragraphs.Item($jupiterIdx).Range.Delete()\n$d.Paragraphs.Item($spacerIdx).Range.Delete()\n"}
